{"js": "// Update the worksheet date and the twenty-five \"three-digit \u00f7 one-digit\"\n// answer cells to the new day's values. Each old value is unique in the\n// document, so a scoped, case-sensitive search-and-replace per pair is\n// unambiguous and safe to run independently of ordering.\nconst replacements = [\n  [\"2024-09-22 Sunday\", \"2024-09-23 Monday\"],\n  [\"694\\u00F76=115, 4\", \"406\\u00F74=101, 2\"],\n  [\"280\\u00F78=35, 0\", \"410\\u00F76=68, 2\"],\n  [\"225\\u00F76=37, 3\", \"978\\u00F78=122, 2\"],\n  [\"519\\u00F79=57, 6\", \"848\\u00F73=282, 2\"],\n  [\"710\\u00F74=177, 2\", \"815\\u00F78=101, 7\"],\n  [\"732\\u00F76=122, 0\", \"331\\u00F77=47, 2\"],\n  [\"852\\u00F74=213, 0\", \"997\\u00F73=332, 1\"],\n  [\"733\\u00F78=91, 5\", \"195\\u00F72=97, 1\"],\n  [\"575\\u00F73=191, 2\", \"962\\u00F75=192, 2\"],\n  [\"425\\u00F78=53, 1\", \"508\\u00F77=72, 4\"],\n  [\"726\\u00F74=181, 2\", \"480\\u00F77=68, 4\"],\n  [\"687\\u00F76=114, 3\", \"831\\u00F76=138, 3\"],\n  [\"952\\u00F74=238, 0\", \"522\\u00F73=174, 0\"],\n  [\"647\\u00F74=161, 3\", \"272\\u00F72=136, 0\"],\n  [\"194\\u00F73=64, 2\", \"906\\u00F72=453, 0\"],\n  [\"315\\u00F74=78, 3\", \"648\\u00F73=216, 0\"],\n  [\"620\\u00F77=88, 4\", \"470\\u00F75=94, 0\"],\n  [\"267\\u00F73=89, 0\", \"872\\u00F78=109, 0\"],\n  [\"757\\u00F74=189, 1\", \"989\\u00F76=164, 5\"],\n  [\"322\\u00F72=161, 0\", \"828\\u00F75=165, 3\"],\n  [\"884\\u00F79=98, 2\", \"543\\u00F75=108, 3\"],\n  [\"349\\u00F77=49, 6\", \"178\\u00F76=29, 4\"],\n  [\"572\\u00F79=63, 5\", \"438\\u00F76=73, 0\"],\n  [\"742\\u00F74=185, 2\", \"725\\u00F73=241, 2\"],\n  [\"488\\u00F79=54, 2\", \"587\\u00F78=73, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the twenty-five \"three-digit \u00f7 one-digit\"\n# answer cells to the new day's values. Each old value is unique in the\n# document, so a document-wide Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2024-09-22 Sunday', '2024-09-23 Monday')\n    ,@('694\u00f76=115, 4', '406\u00f74=101, 2')\n    ,@('280\u00f78=35, 0', '410\u00f76=68, 2')\n    ,@('225\u00f76=37, 3', '978\u00f78=122, 2')\n    ,@('519\u00f79=57, 6', '848\u00f73=282, 2')\n    ,@('710\u00f74=177, 2', '815\u00f78=101, 7')\n    ,@('732\u00f76=122, 0', '331\u00f77=47, 2')\n    ,@('852\u00f74=213, 0', '997\u00f73=332, 1')\n    ,@('733\u00f78=91, 5', '195\u00f72=97, 1')\n    ,@('575\u00f73=191, 2', '962\u00f75=192, 2')\n    ,@('425\u00f78=53, 1', '508\u00f77=72, 4')\n    ,@('726\u00f74=181, 2', '480\u00f77=68, 4')\n    ,@('687\u00f76=114, 3', '831\u00f76=138, 3')\n    ,@('952\u00f74=238, 0', '522\u00f73=174, 0')\n    ,@('647\u00f74=161, 3', '272\u00f72=136, 0')\n    ,@('194\u00f73=64, 2', '906\u00f72=453, 0')\n    ,@('315\u00f74=78, 3', '648\u00f73=216, 0')\n    ,@('620\u00f77=88, 4', '470\u00f75=94, 0')\n    ,@('267\u00f73=89, 0', '872\u00f78=109, 0')\n    ,@('757\u00f74=189, 1', '989\u00f76=164, 5')\n    ,@('322\u00f72=161, 0', '828\u00f75=165, 3')\n    ,@('884\u00f79=98, 2', '543\u00f75=108, 3')\n    ,@('349\u00f77=49, 6', '178\u00f76=29, 4')\n    ,@('572\u00f79=63, 5', '438\u00f76=73, 0')\n    ,@('742\u00f74=185, 2', '725\u00f73=241, 2')\n    ,@('488\u00f79=54, 2', '587\u00f78=73, 3')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
